# Diário de Bordo.docx — apply commit "Autenticação de usuário finalizado
# (backend e frontend)":
#   1. Merge the two runs that make up the "10º dia" heading into one run.
#   2. Append the "11º dia" entry (four new paragraphs) after the GitKraken
#      paragraph, moving the _GoBack bookmark to sit right before the final
#      "." run of the new last paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "10º dia – Dia 16" + "/03/2019"  ->  "10º dia – Dia 16/03/2019"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
[void]$rng.Find.Execute(
    "10º dia – Dia 16/03/2019",  # FindText (already contiguous text; run split doesn't matter to Find)
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "10º dia – Dia 16/03/2019",  # ReplaceWith (same text -> Word rewrites it as a single run)
    2)

# ---------------------------------------------------------------------
# 2) Remove the (hidden) _GoBack bookmark from its old position — it will
#    be re-created later at the new location.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3) Locate the paragraph that ends with "...controle de versionamento."
#    and the trailing two (empty) paragraphs that follow it.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
[void]$rng2.Find.Execute("controle de versionamento.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$gitkrakenPara = $rng2.Paragraphs(1)
$blankPara = $gitkrakenPara.Next()      # first trailing empty paragraph
$lastPara = $blankPara.Next()           # final (sz=20) empty paragraph of the document

# 3a) Turn the first trailing empty paragraph into the spacer paragraph:
#     add spacing-after = 0 (keeps Arial/sz24/jc=both already in place).
$blankPara.SpaceAfter = 0

# ---------------------------------------------------------------------
# 4) Build the four new paragraphs as raw WordprocessingML and insert
#    them at the start of the final (sz=20) paragraph — InsertXML
#    replaces that paragraph's own pPr/run content with the *last*
#    supplied <w:p>, and inserts any preceding <w:p> elements as brand
#    new paragraphs just before it.
# ---------------------------------------------------------------------
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$rPr = "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr>"
$pPr = "<w:pPr><w:jc w:val=`"both`"/>$rPr</w:pPr>"

$para11dia = "<w:p $w>$pPr<w:r>$rPr<w:t>11º dia – 02/04/2019</w:t></w:r></w:p>"

$paraAntes = "<w:p $w>$pPr" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">Antes de ontem, domingo, terminei de configurar a parte de autenticação do sistema. Tanto o backend como o frontend estão funcionando perfeitamente. O cadastro de novo usuário e a criação de uma chave token estão funcionando da forma esperada. Após a autenticacao o usuário é redirecionado para a página inicial da aplicação onde estarão listadas as suas histórias, por </w:t></w:r>" +
    "<w:r>$rPr<w:lastRenderedPageBreak/><w:t>enquanto essa página ainda não foi implementada</w:t></w:r>" +
    "<w:r>$rPr<w:t>, esse será nosso próximo passo.</w:t></w:r>" +
    "</w:p>"

$paraApis = "<w:p $w>$pPr" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">Estamos usando duas api’s no backend, a primeira pode ser acessada por qualquer pessoa, estamos chamando ela de openApi, pois ela é aberta. A openApi dá acesso a três rotas: /login, /signup e /validateToken. A outra api usada é a </w:t></w:r>" +
    "<w:r>$rPr<w:t>protectedApi</w:t></w:r>" +
    "<w:r>$rPr<w:t>, ela só pode ser acessada com um token válido, nela podemos acessar as rotas de cadastro, edição e remoção dos textos narrativos.</w:t></w:r>" +
    "</w:p>"

$paraGitKraken = "<w:p $w>$pPr" +
    "<w:r>$rPr<w:t>Estou usando o GitKraken para gerenciar esse projeto, ele deixa mais intuitivo o versionamento do código</w:t></w:r>" +
    "<w:r>$rPr<w:t>, fiz alguns textes com ele na atividade de criptografia do professor Fábio e gostei do modo do fluxo de trabalho</w:t></w:r>" +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
    "<w:r>$rPr<w:t>.</w:t></w:r>" +
    "</w:p>"

$xml = $para11dia + $paraAntes + $paraApis + $paraGitKraken

$insertRng = $lastPara.Range
$insertRng.Collapse(1)   # wdCollapseStart
[void]$insertRng.InsertXML($xml)
